$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bugs")

# --- New bug report row (row 6) -------------------------------------------
# Column order below mirrors the natural data-entry order an author would
# use: # / Reported Date / Reported By / Severity / Assigned To /
# Short Desc. / Expected / Actual / Steps, then Status last, then the
# (unlabeled) Fixed? column K left blank-but-formatted.

$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = 45628
$ws.Cells.Item(6, 2).NumberFormat = "m/d/yy"
$ws.Cells.Item(6, 3).Value = "Edward"
$ws.Cells.Item(6, 4).Value = "P5"
$ws.Cells.Item(6, 6).Value = "Malcolm"
$ws.Cells.Item(6, 7).Value = "Shooting buckshot then press R may cause some stages to be gone"
$ws.Cells.Item(6, 8).Value = "all statges should be there whenever R is pressed"
$ws.Cells.Item(6, 9).Value = "as per description"
$ws.Cells.Item(6, 10).Value = "as per description"
$ws.Cells.Item(6, 5).Value = "not fixed"

# Column K keeps the wrapped-text formatting used elsewhere in the table
# but is left without a value (matches the "Fixed" column being blank for
# this still-open bug).
$ws.Cells.Item(6, 11).WrapText = $true

# Row heights settle at their new autofit-driven values once the row has
# content.
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 150
$ws.Rows.Item(5).RowHeight = 135
$ws.Rows.Item(6).RowHeight = 60

# The author finished data entry on the Status cell of the new row.
$ws.Range("E6").Select() | Out-Null
